$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text (matches source which
# stores these as literal strings like "0.9989", "1.000", etc.)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.358.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.845.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6363"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07548"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2964"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.843.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.990"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6830"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009903"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.388.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.0000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.372"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.462"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05699"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.246"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.130"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.844"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.157"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7168"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.591"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.257.66"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01806"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9081"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.180"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.001.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.219"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.061"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4024"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1128"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05736"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.24%  "
